# Update cryptocurrency price (D) and 1h volume-change (E) columns
# to reflect the refreshed scrape from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.773.24"
$ws.Range("E2").Value = "  +0.29%  "
# Row 3
$ws.Range("D3").Value = "3.450.54"
$ws.Range("E3").Value = "  -0.21%  "
# Row 4
$ws.Range("E4").Value = "  -0.04%  "
# Row 5
$ws.Range("D5").Value = "'578.85"
$ws.Range("E5").Value = "  +0.32%  "
# Row 6
$ws.Range("D6").Value = "'149.56"
$ws.Range("E6").Value = "  +2.49%  "
# Row 7
$ws.Range("E7").Value = "  -0.08%  "
# Row 8
$ws.Range("E8").Value = "  +1.30%  "
# Row 9
$ws.Range("D9").Value = "'8.02"
$ws.Range("E9").Value = "  +5.40%  "
# Row 10
$ws.Range("E10").Value = "  -0.34%  "
# Row 11
$ws.Range("E11").Value = "  +4.16%  "
# Row 12
$ws.Range("D12").Value = "4.042.44"
$ws.Range("E12").Value = "  -0.14%  "
# Row 13
$ws.Range("E13").Value = "  -0.15%  "
# Row 14
$ws.Range("D14").Value = "'28.28"
$ws.Range("E14").Value = "  -4.72%  "
# Row 15
$ws.Range("D15").Value = "3.444.76"
$ws.Range("E15").Value = "  -0.70%  "
# Row 16
$ws.Range("E16").Value = "  +1.79%  "
# Row 17
$ws.Range("D17").Value = "62.811.50"
$ws.Range("E17").Value = "  +0.17%  "
# Row 18
$ws.Range("D18").Value = "'6.39"
$ws.Range("E18").Value = "  +0.56%  "
# Row 19
$ws.Range("D19").Value = "'14.63"
$ws.Range("E19").Value = "  +1.88%  "
# Row 20
$ws.Range("E20").Value = "  -1.96%  "
# Row 21
$ws.Range("D21").Value = "'387.54"
$ws.Range("E21").Value = "  -0.07%  "
# Row 22
$ws.Range("E22").Value = "  +1.05%  "
# Row 23
$ws.Range("D23").Value = "'75.25"
$ws.Range("E23").Value = "  +0.60%  "
# Row 24
$ws.Range("E24").Value = "  +0.03%  "
# Row 25
$ws.Range("D25").Value = "3.586.15"
$ws.Range("E25").Value = "  -0.50%  "
# Row 26
$ws.Range("E26").Value = "  +1.39%  "
# Row 27
$ws.Range("E27").Value = "  +2.09%  "
# Row 28
$ws.Range("D28").Value = "'7.72"
$ws.Range("E28").Value = "  +2.05%  "
# Row 29
$ws.Range("E29").Value = "  +0.04%  "
# Row 30
$ws.Range("E30").Value = "  -1.02%  "
# Row 31
$ws.Range("E31").Value = "  -0.55%  "
# Row 32
$ws.Range("E32").Value = "  +0.00%  "
# Row 33
$ws.Range("E33").Value = "  -2.87%  "
# Row 34
$ws.Range("D34").Value = "'23.25"
$ws.Range("E34").Value = "  -1.97%  "
# Row 35
$ws.Range("D35").Value = "'5.43"
$ws.Range("E35").Value = "  +3.35%  "
# Row 36
$ws.Range("E36").Value = "  +4.66%  "
# Row 37
$ws.Range("D37").Value = "'32.09"
$ws.Range("E37").Value = "  +3.10%  "
# Row 38
$ws.Range("D38").Value = "'6.95"
$ws.Range("E38").Value = "  -1.67%  "
# Row 39
$ws.Range("D39").Value = "'169.13"
$ws.Range("E39").Value = "  -0.92%  "
# Row 40
$ws.Range("D40").Value = "3.485.04"
$ws.Range("E40").Value = "  -0.36%  "
# Row 41
$ws.Range("E41").Value = "  +1.32%  "
# Row 42
$ws.Range("D42").Value = "'42.88"
$ws.Range("E42").Value = "  +1.75%  "
# Row 43
$ws.Range("E43").Value = "  -1.36%  "
# Row 44
$ws.Range("E44").Value = "  -2.11%  "
# Row 45
$ws.Range("E45").Value = "  -0.69%  "
# Row 46
$ws.Range("E46").Value = "  -0.65%  "
# Row 47
$ws.Range("D47").Value = "2.567.61"
$ws.Range("E47").Value = "  -0.84%  "
# Row 48
$ws.Range("E48").Value = "  +2.25%  "
# Row 49
$ws.Range("E49").Value = "  +3.13%  "
# Row 50
$ws.Range("D50").Value = "'22.64"
$ws.Range("E50").Value = "  -2.72%  "
# Row 51
$ws.Range("E51").Value = "  +0.04%  "
